# aggiornamento fino a 1/09/2021
# Append new daily rows (358-366) to the existing data table, mirroring the
# formatting of the last existing row (357).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(358, 44432, 0, 12, 70.03209804493727),
    @(359, 44433, 2, 14, 81.70411438576014),
    @(360, 44434, 4, 13, 75.8681062153487),
    @(361, 44435, 2, 13, 75.8681062153487),
    @(362, 44436, 1, 13, 75.8681062153487),
    @(363, 44437, 7, 19, 110.8841552378173),
    @(364, 44438, 4, 20, 116.7201634082288),
    @(365, 44439, 10, 30, 175.0802451123432),
    @(366, 44440, 0, 28, 163.4082287715203)
)

foreach ($entry in $newRows) {
    $r = $entry[0]

    # Copy formatting (style/number format) from the last populated row so the
    # new row matches the existing table's look (bold centered date cell, etc.)
    $ws.Range("A357:D357").Copy()
    $ws.Range("A" + $r + ":D" + $r).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
    $ws.Cells.Item($r, 4).Value = $entry[4]
}
